$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'25.770.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.20%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.633.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.16%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  -0.26%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'215.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.06%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = "'  -0.75%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  -0.25%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  +0.12%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  -0.92%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'19.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.47%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.08%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  -0.13%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'1.638.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.02%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'1.858.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.21%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.560"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.02%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("E16").Value = "'  -0.25%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'62.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.76%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'25.780.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.30%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  -0.28%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = "'  +1.73%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'194.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.45%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'9.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.32%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "'  +2.27%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "'  -0.22%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "'  +3.79%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'142.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.07%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "'  -0.39%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  +0.47%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'15.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.01%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  -0.22%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  -0.31%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.12%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  -0.04%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "'  +0.84%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "'  -0.04%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  -0.01%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'1.128.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.47%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  -1.73%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = "'  -1.76%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  -0.48%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  +0.28%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  +2.14%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'100.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.90%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.806"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.74%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'1.768.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.33%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.0₆0108"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.46%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'55.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.88%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "'  -2.41%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.0502"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.12%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'7.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.17%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "'  +2.33%  "
$ws.Range("E51").Style = "Normal"
